$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.606.94"
$ws.Range("E2").Value = "  +4.97%  "
$ws.Range("D3").Value = "3.651.22"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "202.49"
$ws.Range("E5").Value = "  +11.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "572.55"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D7").Value = "3.632.42"
$ws.Range("E7").Value = "  +3.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.620"
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.682"
$ws.Range("E10").Value = "  +3.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.09"
$ws.Range("E11").Value = "  +8.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.155"
$ws.Range("E12").Value = "  +8.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000295"
$ws.Range("E13").Value = "  +16.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.10"
$ws.Range("E14").Value = "  +4.14%  "
$ws.Range("D15").Value = "4.216.51"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("D16").Value = "3.643.37"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "68.433.84"
$ws.Range("E18").Value = "  +4.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.71"
$ws.Range("E19").Value = "  +3.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.47"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("E21").Value = "  +4.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "403.77"
$ws.Range("E22").Value = "  +3.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.23"
$ws.Range("E23").Value = "  +28.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.23"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.91"
$ws.Range("E25").Value = "  +2.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.97"
$ws.Range("E26").Value = "  +4.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.60"
$ws.Range("E27").Value = "  +2.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.86"
$ws.Range("E28").Value = "  +9.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.13"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.20"
$ws.Range("E30").Value = "  +21.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.19"
$ws.Range("E31").Value = "  +4.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.01"
$ws.Range("E32").Value = "  +4.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "697.34"
$ws.Range("E33").Value = "  +15.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.29"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.117"
$ws.Range("E35").Value = "  +5.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.43"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.90"
$ws.Range("E37").Value = "  +5.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.425"
$ws.Range("E38").Value = "  +15.38%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "0.0₃0786"
$ws.Range("E40").Value = "  +6.07%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.280.62"
$ws.Range("E41").Value = "  +15.36%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.139"
$ws.Range("E42").Value = "  +7.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.15"
$ws.Range("E43").Value = "  +14.33%  "
$ws.Range("E44").Value = "  +17.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.05"
$ws.Range("E45").Value = "  +37.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.995"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0421"
$ws.Range("E47").Value = "  +4.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  +12.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.98"
$ws.Range("E49").Value = "  +9.98%  "
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("E51").Value = "  +6.38%  "
